$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F
$ws.Range("F1").Value = "test column"

# Repeating pattern test_a, test_b, test_c for rows 2-13
$values = @("test_a", "test_b", "test_c")
for ($r = 2; $r -le 13; $r++) {
    $idx = ($r - 2) % 3
    $ws.Cells.Item($r, 6).Value = $values[$idx]
}

# Update the selection to match the authored state (F17)
$ws.Range("F17").Select()
